$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.2147957804815247, 0.3580819848962541, 0.2088379301833887, 0.4569878884427777, 0.4185882721160311, 14)
    3  = @(0.245473867636748,  0.3683843928281823, 0.2358514199634523, 0.4856453644002507, 0.4361500127557235, 13)
    4  = @(0.301463903713677,  0.3567050099425655, 0.2583510523008795, 0.5082824532687308, 0.4274285484895888, 12)
    5  = @(0.2900334324511264, 0.3439303588265401, 0.1874043702602189, 0.4329022640968963, 0.3370659818980856, 11)
    6  = @(0.2833661826906564, 0.2909807410604855, 0.2099604875429734, 0.4582144558424291, 0.3795671424286398, 10)
    7  = @(0.3070978867771534, 0.3733636115588901, 0.2272551859511885, 0.476712896774556,  0.3867354821117653, 9)
    8  = @(0.3125739411747067, 0.4101416856701894, 0.3161064439746084, 0.5622334425971194, 0.499604386659364,  8)
    9  = @(0.2627214152773011, 0.3235434349789264, 0.1669924024756296, 0.4086470389904099, 0.3380801243818421, 7)
    10 = @(0.3023013051968123, 0.3856149542229345, 0.3827274330185637, 0.618649685216572,  0.5912779588939882, 6)
    11 = @(0.302435154479656,  0.4668197381621303, 0.3675323742590614, 0.6062444838999044, 0.5874365408215683, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 6).Value = $vals[4]  # F
    $ws.Cells.Item($row, 7).Value = $vals[5]  # G
}

$wb.Save()
